# Update simulation results in Sheet1 for the "380 kV" case.
# Rows 2-25 (A column holds index 0..23) contain per-row metrics in columns
# B..O (columns G, I, J, K, N are unchanged zeros). This applies the updated
# values as produced by the new case run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.118015681109569
$ws.Range("C2").Value = 0.3082411197244141
$ws.Range("D2").Value = 0.03355238475187861
$ws.Range("E2").Value = 0.1233058829017842
$ws.Range("F2").Value = 0.8231263785664211
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("L2").Value = 0.1891530951563141
$ws.Range("M2").Value = 0.2361424816765947
$ws.Range("O2").Value = 2.871797905458948
$ws.Range("B3").Value = 1.012251070564901
$ws.Range("C3").Value = 0.2946074768429412
$ws.Range("D3").Value = 0.03206211263224645
$ws.Range("E3").Value = 0.124723630836637
$ws.Range("F3").Value = 0.8238021451527047
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("L3").Value = 0.1865257311867481
$ws.Range("M3").Value = 0.2199033853709622
$ws.Range("O3").Value = 2.890035906224853
$ws.Range("B4").Value = 0.9473682001214456
$ws.Range("C4").Value = 0.2862065952481032
$ws.Range("D4").Value = 0.0311432607357176
$ws.Range("E4").Value = 0.125643821480906
$ws.Range("F4").Value = 0.8248677082826745
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("L4").Value = 0.1850045961093016
$ws.Range("M4").Value = 0.2099848639771196
$ws.Range("O4").Value = 2.903466340057037
$ws.Range("B5").Value = 0.9209436966755788
$ws.Range("C5").Value = 0.2827759353475017
$ws.Range("D5").Value = 0.03076788752294846
$ws.Range("E5").Value = 0.1260313135664947
$ws.Range("F5").Value = 0.8254652868235866
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("L5").Value = 0.1844079178932034
$ws.Range("M5").Value = 0.2059563668816935
$ws.Range("O5").Value = 2.90949973733305
$ws.Range("B6").Value = 0.9165569249604459
$ws.Range("C6").Value = 0.2822058475786093
$ws.Range("D6").Value = 0.03070550149033835
$ws.Range("E6").Value = 0.1260964120891777
$ws.Range("F6").Value = 0.8255743728934988
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("L6").Value = 0.1843102422320513
$ws.Range("M6").Value = 0.2052882526239443
$ws.Range("O6").Value = 2.9105353971635
$ws.Range("B7").Value = 0.9470117641438947
$ws.Range("C7").Value = 0.2861603571142837
$ws.Range("D7").Value = 0.03113820206615259
$ws.Range("E7").Value = 0.1256489966780163
$ws.Range("F7").Value = 0.824875106367891
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("L7").Value = 0.1849964551176484
$ws.Range("M7").Value = 0.209930479787964
$ws.Range("O7").Value = 2.903545440816259
$ws.Range("B8").Value = 1.081537041416539
$ws.Range("C8").Value = 0.3035465653591132
$ws.Range("D8").Value = 0.03303934696928224
$ws.Range("E8").Value = 0.123784417381288
$ws.Range("F8").Value = 0.8232241865317249
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("L8").Value = 0.1882280939913983
$ws.Range("M8").Value = 0.2305325005613525
$ws.Range("O8").Value = 2.877622509305354
$ws.Range("B9").Value = 1.345743109301338
$ws.Range("C9").Value = 0.3373952499985364
$ws.Range("D9").Value = 0.03673620094449603
$ws.Range("E9").Value = 0.1205216274595068
$ws.Range("F9").Value = 0.8251620350377564
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("L9").Value = 0.1952945804903692
$ws.Range("M9").Value = 0.2713410357480726
$ws.Range("O9").Value = 2.844543455190006
$ws.Range("B10").Value = 1.540053577822903
$ws.Range("C10").Value = 0.3621040387105836
$ws.Range("D10").Value = 0.03943215737670869
$ws.Range("E10").Value = 0.118363539648797
$ws.Range("F10").Value = 0.8297602121856897
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("L10").Value = 0.2009299290899378
$ws.Range("M10").Value = 0.3015652178322341
$ws.Range("O10").Value = 2.831128854499269
$ws.Range("B11").Value = 1.628484986013518
$ws.Range("C11").Value = 0.3733080582367734
$ws.Range("D11").Value = 0.04065404735843003
$ws.Range("E11").Value = 0.117433493432773
$ws.Range("F11").Value = 0.8325456133071754
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("L11").Value = 0.203589771150817
$ws.Range("M11").Value = 0.3153663956591686
$ws.Range("O11").Value = 2.827404545918114
$ws.Range("B12").Value = 1.661976073157291
$ws.Range("C12").Value = 0.3775453212936384
$ws.Range("D12").Value = 0.04111607394087713
$ws.Range("E12").Value = 0.1170887276306282
$ws.Range("F12").Value = 0.8337004295917723
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("L12").Value = 0.2046108018384984
$ws.Range("M12").Value = 0.3205998604776923
$ws.Range("O12").Value = 2.826337285004115
$ws.Range("B13").Value = 1.654763012999013
$ws.Range("C13").Value = 0.3766329976444922
$ws.Range("D13").Value = 0.04101659877743913
$ws.Range("E13").Value = 0.1171626491358583
$ws.Range("F13").Value = 0.8334472649934241
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("L13").Value = 0.2043902912210456
$ws.Range("M13").Value = 0.3194724209483653
$ws.Range("O13").Value = 2.826551863985969
$ws.Range("B14").Value = 1.631240245745175
$ws.Range("C14").Value = 0.3736567712431906
$ws.Range("D14").Value = 0.04069207227984606
$ws.Range("E14").Value = 0.1174049806346802
$ws.Range("F14").Value = 0.8326386139872142
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("L14").Value = 0.2036734954947121
$ws.Range("M14").Value = 0.3157968117782346
$ws.Range("O14").Value = 2.827309860172193
$ws.Range("B15").Value = 1.616832351879168
$ws.Range("C15").Value = 0.3718330284405624
$ws.Range("D15").Value = 0.04049320150195257
$ws.Range("E15").Value = 0.1175543820454277
$ws.Range("F15").Value = 0.8321563295279617
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("L15").Value = 0.2032362341902569
$ws.Range("M15").Value = 0.3135463336099491
$ws.Range("O15").Value = 2.827818862790622
$ws.Range("B16").Value = 1.534275023905366
$ws.Range("C16").Value = 0.3613710802930825
$ws.Range("D16").Value = 0.039352210797297
$ws.Range("E16").Value = 0.1184253595077182
$ws.Range("F16").Value = 0.8295921634482113
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("L16").Value = 0.2007580363463575
$ws.Range("M16").Value = 0.3006643063454746
$ws.Range("O16").Value = 2.831420192909405
$ws.Range("B17").Value = 1.483637578443791
$ws.Range("C17").Value = 0.3549435734134079
$ws.Range("D17").Value = 0.03865107409946944
$ws.Range("E17").Value = 0.1189729058223672
$ws.Range("F17").Value = 0.828197007259007
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("L17").Value = 0.1992623763262031
$ws.Range("M17").Value = 0.2927747686860585
$ws.Range("O17").Value = 2.834239348906095
$ws.Range("B18").Value = 1.454515982280839
$ws.Range("C18").Value = 0.3512432558491412
$ws.Range("D18").Value = 0.03824737554735691
$ws.Range("E18").Value = 0.1192927047498393
$ws.Range("F18").Value = 0.8274598239411262
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("L18").Value = 0.1984111782779934
$ws.Range("M18").Value = 0.2882418274550957
$ws.Range("O18").Value = 2.836084622194306
$ws.Range("B19").Value = 1.444656596184245
$ws.Range("C19").Value = 0.3499898184026904
$ws.Range("D19").Value = 0.03811061830187157
$ws.Range("E19").Value = 0.1194018189715029
$ws.Range("F19").Value = 0.8272214280706862
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("L19").Value = 0.1981245355892298
$ws.Range("M19").Value = 0.2867079020250927
$ws.Range("O19").Value = 2.836747799293249
$ws.Range("B20").Value = 1.4890276466773
$ws.Range("C20").Value = 0.3556281447869765
$ws.Range("D20").Value = 0.0387257553048741
$ws.Range("E20").Value = 0.1189141151613362
$ws.Range("F20").Value = 0.8283387662841619
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("L20").Value = 0.199420653883962
$ws.Range("M20").Value = 0.2936141166913231
$ws.Range("O20").Value = 2.833916077721454
$ws.Range("B21").Value = 1.638149355285805
$ws.Range("C21").Value = 0.3745311112121783
$ws.Range("D21").Value = 0.04078741216056159
$ws.Range("E21").Value = 0.117333600624101
$ws.Range("F21").Value = 0.8328734168322001
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("L21").Value = 0.203883661401008
$ws.Range("M21").Value = 0.316876232090479
$ws.Range("O21").Value = 2.827077899255642
$ws.Range("B22").Value = 1.735631820369292
$ws.Range("C22").Value = 0.3868533724960344
$ws.Range("D22").Value = 0.04213086761872376
$ws.Range("E22").Value = 0.11634390309103
$ws.Range("F22").Value = 0.836420350245362
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("L22").Value = 0.206880949371822
$ws.Range("M22").Value = 0.3321215728849793
$ws.Range("O22").Value = 2.824608838518543
$ws.Range("B23").Value = 1.683602005027296
$ws.Range("C23").Value = 0.3802797555388508
$ws.Range("D23").Value = 0.04141421117680721
$ws.Range("E23").Value = 0.1168681678770462
$ws.Range("F23").Value = 0.8344738176604238
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("L23").Value = 0.2052738921474884
$ws.Range("M23").Value = 0.3239810600088546
$ws.Range("O23").Value = 2.825743255344719
$ws.Range("B24").Value = 1.486590826167173
$ws.Range("C24").Value = 0.3553186658066068
$ws.Range("D24").Value = 0.03869199381659882
$ws.Range("E24").Value = 0.1189406788213578
$ws.Range("F24").Value = 0.8282744748573094
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("L24").Value = 0.1993490695703741
$ws.Range("M24").Value = 0.2932346385362479
$ws.Range("O24").Value = 2.834061529383575
$ws.Range("B25").Value = 1.274230211854615
$ws.Range("C25").Value = 0.3282656365030334
$ws.Range("D25").Value = 0.0357395631526316
$ws.Range("E25").Value = 0.1213622489075763
$ws.Range("F25").Value = 0.824081728257994
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("L25").Value = 0.1933049387016226
$ws.Range("M25").Value = 0.2602582533825384
$ws.Range("O25").Value = 2.851584689225348
